$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 60.90318633333334
$ws.Range("H2").Value = 182.709559
$ws.Range("I2").Value = 0.4799022665420342
$ws.Range("J2").Value = 0.4799022665420342
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.481206
$ws.Range("N2").Value = 262.443618
$ws.Range("O2").Value = 0.8890954013316028
$ws.Range("P2").Value = 0.8890954013316029
$ws.Range("Q2").Value = 5327.884189682718
$ws.Range("R2").Value = 47950.95770714447
$ws.Range("S2").Value = 0.4266788982711357
$ws.Range("T2").Value = 0.4266788982711358

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 60.90318633333334
$ws.Range("H3").Value = 182.709559
$ws.Range("I3").Value = 0.4799022665420342
$ws.Range("J3").Value = 0.4799022665420342
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.36557633333333
$ws.Range("N3").Value = 31.096729
$ws.Range("O3").Value = 0.1053481847303107
$ws.Range("P3").Value = 0.1053481847303107
$ws.Range("Q3").Value = 631.2966268813901
$ws.Range("R3").Value = 5681.669641932512
$ws.Range("S3").Value = 0.05055683262816501
$ws.Range("T3").Value = 0.05055683262816502

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 60.90318633333334
$ws.Range("H4").Value = 182.709559
$ws.Range("I4").Value = 0.4799022665420342
$ws.Range("J4").Value = 0.4799022665420342
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.546715
$ws.Range("N4").Value = 1.640145
$ws.Range("O4").Value = 0.005556413938086396
$ws.Range("P4").Value = 0.005556413938086396
$ws.Range("Q4").Value = 33.29668551622833
$ws.Range("R4").Value = 299.670169646055
$ws.Range("S4").Value = 0.002666535642733411
$ws.Range("T4").Value = 0.002666535642733411

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.42568199999999
$ws.Range("H5").Value = 154.277046
$ws.Range("I5").Value = 0.405221842009972
$ws.Range("J5").Value = 0.405221842009972
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 87.481206
$ws.Range("N5").Value = 262.443618
$ws.Range("O5").Value = 0.8890954013316028
$ws.Range("P5").Value = 0.8890954013316029
$ws.Range("Q5").Value = 4498.780680732491
$ws.Range("R5").Value = 40489.02612659243
$ws.Range("S5").Value = 0.3602808762501874
$ws.Range("T5").Value = 0.3602808762501874

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 51.42568199999999
$ws.Range("H6").Value = 154.277046
$ws.Range("I6").Value = 0.405221842009972
$ws.Range("J6").Value = 0.405221842009972
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.36557633333333
$ws.Range("N6").Value = 31.096729
$ws.Range("O6").Value = 0.1053481847303107
$ws.Range("P6").Value = 0.1053481847303107
$ws.Range("Q6").Value = 533.0568322647259
$ws.Range("R6").Value = 4797.511490382533
$ws.Range("S6").Value = 0.0426893854688233
$ws.Range("T6").Value = 0.0426893854688233

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 51.42568199999999
$ws.Range("H7").Value = 154.277046
$ws.Range("I7").Value = 0.405221842009972
$ws.Range("J7").Value = 0.405221842009972
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.546715
$ws.Range("N7").Value = 1.640145
$ws.Range("O7").Value = 0.005556413938086396
$ws.Range("P7").Value = 0.005556413938086396
$ws.Range("Q7").Value = 28.11519173463
$ws.Range("R7").Value = 253.03672561167
$ws.Range("S7").Value = 0.002251580290961252
$ws.Range("T7").Value = 0.002251580290961252

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.57860966666667
$ws.Range("H8").Value = 43.735829
$ws.Range("I8").Value = 0.1148758914479938
$ws.Range("J8").Value = 0.1148758914479938
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 87.481206
$ws.Range("N8").Value = 262.443618
$ws.Range("O8").Value = 0.8890954013316028
$ws.Range("P8").Value = 0.8890954013316029
$ws.Range("Q8").Value = 1275.354355443258
$ws.Range("R8").Value = 11478.18919898932
$ws.Range("S8").Value = 0.1021356268102797
$ws.Range("T8").Value = 0.1021356268102797

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.57860966666667
$ws.Range("H9").Value = 43.735829
$ws.Range("I9").Value = 0.1148758914479938
$ws.Range("J9").Value = 0.1148758914479938
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.36557633333333
$ws.Range("N9").Value = 31.096729
$ws.Range("O9").Value = 0.1053481847303107
$ws.Range("P9").Value = 0.1053481847303107
$ws.Range("Q9").Value = 151.1156913337045
$ws.Range("R9").Value = 1360.041222003341
$ws.Range("S9").Value = 0.01210196663332237
$ws.Range("T9").Value = 0.01210196663332237

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.57860966666667
$ws.Range("H10").Value = 43.735829
$ws.Range("I10").Value = 0.1148758914479938
$ws.Range("J10").Value = 0.1148758914479938
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.546715
$ws.Range("N10").Value = 1.640145
$ws.Range("O10").Value = 0.005556413938086396
$ws.Range("P10").Value = 0.005556413938086396
$ws.Range("Q10").Value = 7.970344583911665
$ws.Range("R10").Value = 71.73310125520499
$ws.Range("S10").Value = 0.0006382980043917327
$ws.Range("T10").Value = 0.0006382980043917327

Write-Host "Updated values successfully"
